$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the date-cell formatting from the row above (C9) into C10 so the new
# date picks up the same style (fill + date number format) used by the rest
# of the "odd" data rows.
$ws.Range("C9").Copy()
$ws.Range("C10").PasteSpecial(-4122)
$ws.Range("C10").Value = 44027

# Fill in the rest of row 10's data (Día 8 of the log).
$ws.Range("D10").Value = 7
$ws.Range("E10").Value = 'Creacion de diseño de vistas para par el programa "Pantallas"'

# Update the selected cell, matching the author's last selection.
$ws.Range("F16").Select()
